# Update the question texts that received a clarifying parenthetical remark.
# (These are the same Hebrew survey questions, just reworded with examples.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A23").Value = "מאז הסקר האחרון, היה לי קשה להפסיק לעשות משהו אחרי שביקשו ממני להפסיק (לדוגמה כשמורים מבקשים ממני להפסיק לדבר בכיתה או כשאבא מבקש ממני להפסיק להיות מול המסך)"
$ws.Range("A27").Value = "היום אבא שלי עזר לי במשהו (כמו הקפצה לחוג/תנועת נוער, עזרה בשיעורי בית וכו')"
$ws.Range("A28").Value = "היום אמא שלי עזרה לי במשהו (כמו הקפצה לחוג/תנועת נוער, עזרה בשיעורי בית וכו')"
$ws.Range("A44").Value = "מאז הסקר האחרון, שיתפתי את אבא שלי ברגשות/תחושות שלי (שאני שמח / עצוב / כועס / עצבני / שקרה לי משהו מרגש היום)"
$ws.Range("A45").Value = "מאז הסקר האחרון, שיתפתי את אמא שלי ברגשות/תחושות שלי (שאני שמח / עצוב / כועס / עצבני / שקרה לי משהו מרגש היום)"

# Restore the sheet's view/selection state (scrolled down, B44 selected)
# as left by the author after editing the last row that changed.
$excel.Goto($ws.Range("B44"), $true)
$ws.Range("B44").Select()
